$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44348
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("P2").Value = 194

# Row 3
$ws.Range("D3").Value = 44354
$ws.Range("J3").Value = 150
$ws.Range("O3").Value = 'Región del Maule'

# Row 4
$ws.Range("D4").Value = 44342

# Row 5
$ws.Range("D5").Value = 44371
$ws.Range("K5").Value = 6500
$ws.Range("L5").Value = 6500
$ws.Range("M5").Value = 6500
$ws.Range("P5").Value = 181

# Row 7
$ws.Range("D7").Value = 44376
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = 6500
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 181

# Row 8
$ws.Range("D8").Value = 44364
$ws.Range("J8").Value = 100
$ws.Range("O8").Value = 'Región Metropolitana'

# Row 9
$ws.Range("D9").Value = 44355
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("P9").Value = 194

# Row 10
$ws.Range("D10").Value = 44358
$ws.Range("J10").Value = 150
$ws.Range("N10").Value = '$/caja 36 atados'
$ws.Range("P10").Value = 194
$ws.Range("Q10").Value = 36

# Row 11
$ws.Range("D11").Value = 44362
$ws.Range("J11").Value = 100

# Row 13
$ws.Range("D13").Value = 44386
$ws.Range("J13").Value = 200
$ws.Range("N13").Value = '$/caja 36 atados'
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("P13").Value = 181
$ws.Range("Q13").Value = 36

# Row 14
$ws.Range("D14").Value = 44369
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("N14").Value = '$/caja 20 docenas'
$ws.Range("P14").Value = 7000
$ws.Range("Q14").Value = 1

# Row 15
$ws.Range("D15").Value = 44357
$ws.Range("K15").Value = 6500
$ws.Range("L15").Value = 6500
$ws.Range("M15").Value = 6500
$ws.Range("N15").Value = '$/caja 20 docenas'
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 6500
$ws.Range("Q15").Value = 1
